$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the thrust-test data table (columns A:D, rows 2-6).
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 1000
$ws.Cells.Item(2, 3).Value = 90
$ws.Cells.Item(2, 4).Value = 45

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = 1000
$ws.Cells.Item(3, 3).Value = 90
$ws.Cells.Item(3, 4).Value = 45

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 1000
$ws.Cells.Item(4, 3).Value = 90
$ws.Cells.Item(4, 4).Value = 45

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = 900
$ws.Cells.Item(5, 3).Value = 80
$ws.Cells.Item(5, 4).Value = 45

$ws.Cells.Item(6, 1).Value = 10
$ws.Cells.Item(6, 2).Value = 1000
$ws.Cells.Item(6, 3).Value = 60
$ws.Cells.Item(6, 4).Value = 45

# Move the active selection to match the saved view state.
$ws.Range("G17").Select()
